$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 805.6667
$ws.Range("J4").Value = 1500.5
$ws.Range("L4").Value = 1500.5
$ws.Range("N4").Value = -1728.5
$ws.Range("H15").Value = 1017.9889
$ws.Range("I15").Value = 1017.9889
$ws.Range("K15").Value = 3053.9667
$ws.Range("M15").Value = -2884.9667
$ws.Range("H18").Value = 681.5714
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H28").Value = 3082.5483
$ws.Range("I28").Value = 3040.2415
$ws.Range("K28").Value = 3040.2415
$ws.Range("M28").Value = -2555.2415
$ws.Range("H43").Value = 1392.8572
$ws.Range("I43").Value = 720.3333
$ws.Range("K43").Value = 720.3333
$ws.Range("M43").Value = -651.3333
$ws.Range("H53").Value = 1286.75
$ws.Range("I53").Value = 2228.2
$ws.Range("J53").Value = 858.8182
$ws.Range("K53").Value = 2228.2
$ws.Range("L53").Value = 858.8182
$ws.Range("M53").Value = -1591.2
$ws.Range("N53").Value = -2132.8182
$ws.Range("H103").Value = 298.42856
$ws.Range("J103").Value = 275
$ws.Range("L103").Value = 825
$ws.Range("N103").Value = -1997
$ws.Range("H113").Value = 6170.391
$ws.Range("I113").Value = 6042
$ws.Range("K113").Value = 6042
$ws.Range("M113").Value = -2788
$ws.Range("H125").Value = 5720.5557
$ws.Range("I125").Value = 4889.8
$ws.Range("J125").Value = 6759
$ws.Range("K125").Value = 44008.2
$ws.Range("L125").Value = 60831
$ws.Range("M125").Value = -41548.2
$ws.Range("N125").Value = -65751
$ws.Range("H134").Value = 34866.938
$ws.Range("J134").Value = 34866.938
$ws.Range("L134").Value = 34866.938
$ws.Range("N134").Value = -45006.938
$ws.Range("H138").Value = 3199.41
$ws.Range("I138").Value = 2471.8
$ws.Range("J138").Value = 3327.8118
$ws.Range("K138").Value = 7415.400000000001
$ws.Range("L138").Value = 9983.4354
$ws.Range("M138").Value = -2275.400000000001
$ws.Range("N138").Value = -20263.4354
$ws.Range("H141").Value = 7561.077
$ws.Range("I141").Value = 5046.8125
$ws.Range("K141").Value = 15140.4375
$ws.Range("M141").Value = -9960.4375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1068.85
$ws.Range("I2").Value = 941.94116
$ws.Range("J2").Value = 1788
$ws.Range("K2").Value = 941.94116
$ws.Range("L2").Value = 1788
$ws.Range("M2").Value = -828.94116
$ws.Range("N2").Value = -2014
$ws.Range("H32").Value = 7743.9585
$ws.Range("I32").Value = 4558.3066
$ws.Range("J32").Value = 27495
$ws.Range("K32").Value = 4558.3066
$ws.Range("L32").Value = 27495
$ws.Range("M32").Value = -4271.3066
$ws.Range("N32").Value = -28069
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()
$ws.Range("H74").Value = 3314.05
$ws.Range("I74").Value = 1620.0769
$ws.Range("J74").Value = 6460
$ws.Range("K74").Value = 1620.0769
$ws.Range("L74").Value = 6460
$ws.Range("M74").Value = -746.0769
$ws.Range("N74").Value = -8208
$ws.Range("H77").Value = 3314.05
$ws.Range("I77").Value = 1620.0769
$ws.Range("J77").Value = 6460
$ws.Range("K77").Value = 8100.3845
$ws.Range("L77").Value = 32300
$ws.Range("M77").Value = -3732.3845
$ws.Range("N77").Value = -41036
$ws.Range("H108").Value = 99000
$ws.Range("J108").Value = 99000
$ws.Range("L108").Value = 99000
$ws.Range("N108").Value = -106680
$ws.Range("H116").Value = 1068.85
$ws.Range("I116").Value = 941.94116
$ws.Range("J116").Value = 1788
$ws.Range("K116").Value = 941.94116
$ws.Range("L116").Value = 1788
$ws.Range("M116").Value = 1352.05884
$ws.Range("N116").Value = -6376
$ws.Range("H122").Value = 22998.857
$ws.Range("I122").Value = 62997.5
$ws.Range("K122").Value = 188992.5
$ws.Range("M122").Value = -186542.5
$ws.Range("H124").Value = 119237.664
$ws.Range("J124").Value = 119237.664
$ws.Range("L124").Value = 119237.664
$ws.Range("N124").Value = -129057.664
$ws.Range("H125").Value = 137323.72
$ws.Range("J125").Value = 137323.72
$ws.Range("L125").Value = 137323.72
$ws.Range("N125").Value = -147163.72
$ws.Range("H131").Value = 64999.5
$ws.Range("J131").Value = 64999.5
$ws.Range("L131").Value = 64999.5
$ws.Range("N131").Value = -75079.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1068.85
$ws.Range("I3").Value = 941.94116
$ws.Range("J3").Value = 1788
$ws.Range("K3").Value = 941.94116
$ws.Range("L3").Value = 1788
$ws.Range("M3").Value = -827.94116
$ws.Range("N3").Value = -2016
$ws.Range("H132").Value = 97999.5
$ws.Range("J132").Value = 97999.5
$ws.Range("L132").Value = 97999.5
$ws.Range("N132").Value = -108119.5
$ws.Range("H133").Value = 70993
$ws.Range("J133").Value = 70993
$ws.Range("L133").Value = 70993
$ws.Range("N133").Value = -81113

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 56157.26
$ws.Range("I31").Value = 3505.9
$ws.Range("J31").Value = 114658.78
$ws.Range("K31").Value = 3505.9
$ws.Range("L31").Value = 114658.78
$ws.Range("M31").Value = -3210.9
$ws.Range("N31").Value = -115248.78
$ws.Range("H34").Value = 56157.26
$ws.Range("I34").Value = 3505.9
$ws.Range("J34").Value = 114658.78
$ws.Range("K34").Value = 3505.9
$ws.Range("L34").Value = 114658.78
$ws.Range("M34").Value = -3303.9
$ws.Range("N34").Value = -115062.78
$ws.Range("H99").Value = 3014.3
$ws.Range("I99").Value = 2957.8333
$ws.Range("J99").Value = 3099
$ws.Range("K99").Value = 2957.8333
$ws.Range("L99").Value = 3099
$ws.Range("M99").Value = -1459.8333
$ws.Range("N99").Value = -6095
$ws.Range("H126").Value = 3014.3
$ws.Range("I126").Value = 2957.8333
$ws.Range("J126").Value = 3099
$ws.Range("K126").Value = 8873.499899999999
$ws.Range("L126").Value = 9297
$ws.Range("M126").Value = -6403.499899999999
$ws.Range("N126").Value = -14237
$ws.Range("H132").Value = 3852.05
$ws.Range("I132").Value = 3787.9443
$ws.Range("J132").Value = 4429
$ws.Range("K132").Value = 11363.8329
$ws.Range("L132").Value = 13287
$ws.Range("M132").Value = -8833.832900000001
$ws.Range("N132").Value = -18347

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 168.14285
$ws.Range("I2").Value = 134.5
$ws.Range("J2").Value = 213
$ws.Range("K2").Value = 807
$ws.Range("L2").Value = 1278
$ws.Range("M2").Value = -694
$ws.Range("N2").Value = -1504

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 334.25925
$ws.Range("J2").Value = 794.7
$ws.Range("L2").Value = 794.7
$ws.Range("N2").Value = -1020.7
$ws.Range("H102").Value = 10821.036
$ws.Range("I102").Value = 11600.04
$ws.Range("K102").Value = 11600.04
$ws.Range("M102").Value = -9978.040000000001
$ws.Range("H132").Value = 5605.933
$ws.Range("I132").Value = 4379.6
$ws.Range("J132").Value = 6832.2666
$ws.Range("K132").Value = 13138.8
$ws.Range("L132").Value = 20496.7998
$ws.Range("M132").Value = -10608.8
$ws.Range("N132").Value = -25556.7998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 5786.6
$ws.Range("I4").Value = 5786.6
$ws.Range("K4").Value = 5786.6
$ws.Range("M4").Value = -5673.6
$ws.Range("H28").Value = 5786.6
$ws.Range("I28").Value = 5786.6
$ws.Range("K28").Value = 5786.6
$ws.Range("M28").Value = -5554.6
$ws.Range("H37").Value = 5786.6
$ws.Range("I37").Value = 5786.6
$ws.Range("K37").Value = 5786.6
$ws.Range("M37").Value = -5679.6
$ws.Range("H122").Value = 5553.5557
$ws.Range("I122").Value = 5407.6
$ws.Range("K122").Value = 16222.8
$ws.Range("M122").Value = -13772.8
$ws.Range("H132").Value = 3595.6592
$ws.Range("I132").Value = 3291.5588
$ws.Range("K132").Value = 9874.6764
$ws.Range("M132").Value = -7344.6764

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 34998.75
$ws.Range("I24").Value = 29997.5
$ws.Range("J24").Value = 40000
$ws.Range("K24").Value = 29997.5
$ws.Range("L24").Value = 40000
$ws.Range("M24").Value = -29767.5
$ws.Range("N24").Value = -40460
$ws.Range("H62").Value = 6221
$ws.Range("I62").Value = 4748.75
$ws.Range("K62").Value = 4748.75
$ws.Range("M62").Value = -4124.75
$ws.Range("H65").Value = 6221
$ws.Range("I65").Value = 4748.75
$ws.Range("K65").Value = 23743.75
$ws.Range("M65").Value = -20623.75
$ws.Range("H107").Value = 1264.1
$ws.Range("I107").Value = 1369.8572
$ws.Range("J107").Value = 1207.1538
$ws.Range("K107").Value = 4109.571599999999
$ws.Range("L107").Value = 3621.4614
$ws.Range("M107").Value = -2189.571599999999
$ws.Range("N107").Value = -7461.4614
$ws.Range("H123").Value = 40000
$ws.Range("J123").Value = 40000
$ws.Range("L123").Value = 40000
$ws.Range("N123").Value = -49800
$ws.Range("H131").Value = 80000
$ws.Range("J131").Value = 80000
$ws.Range("L131").Value = 80000
$ws.Range("N131").Value = -90080
$ws.Range("H132").Value = 3399.4285
$ws.Range("I132").Value = 3207.28
$ws.Range("J132").Value = 3879.8
$ws.Range("K132").Value = 9621.84
$ws.Range("L132").Value = 11639.4
$ws.Range("M132").Value = -7091.84
$ws.Range("N132").Value = -16699.4
